$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D-column values stay as text (avoid Excel auto-converting numeric-looking
# strings like "542.80" or "1.00" into actual numbers), while keeping the cell
# style identical to the original (no explicit style index).

$dRows = @(2,3,5,6,8,9,10,11,12,13,15,16,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,36,37,38,39,40,42,43,44,46,48,49,51)
foreach ($r in $dRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

$ws.Range("D2").Value = '61.724.82'
$ws.Range("E2").Value = '  -2.03%  '
$ws.Range("D3").Value = '3.008.05'
$ws.Range("E3").Value = '  -2.04%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '542.80'
$ws.Range("E5").Value = '  +0.42%  '
$ws.Range("D6").Value = '132.32'
$ws.Range("E6").Value = '  -3.85%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Value = '3.002.45'
$ws.Range("E8").Value = '  -1.99%  '
$ws.Range("D9").Value = '0.490'
$ws.Range("E9").Value = '  -0.64%  '
$ws.Range("D10").Value = '6.10'
$ws.Range("E10").Value = '  -2.56%  '
$ws.Range("D11").Value = '0.146'
$ws.Range("E11").Value = '  -6.55%  '
$ws.Range("D12").Value = '0.446'
$ws.Range("E12").Value = '  -1.96%  '
$ws.Range("D13").Value = '34.45'
$ws.Range("E13").Value = '  -0.29%  '
$ws.Range("E14").Value = '  -1.69%  '
$ws.Range("D15").Value = '3.493.45'
$ws.Range("E15").Value = '  -2.11%  '
$ws.Range("D16").Value = '61.784.23'
$ws.Range("E16").Value = '  -1.96%  '
$ws.Range("E17").Value = '  -2.79%  '
$ws.Range("D18").Value = '3.006.75'
$ws.Range("E18").Value = '  -2.18%  '
$ws.Range("D19").Value = '6.62'
$ws.Range("E19").Value = '  -0.24%  '
$ws.Range("D20").Value = '484.06'
$ws.Range("E20").Value = '  +3.06%  '
$ws.Range("D21").Value = '13.23'
$ws.Range("E21").Value = '  -2.54%  '
$ws.Range("D22").Value = '0.668'
$ws.Range("E22").Value = '  -3.98%  '
$ws.Range("D23").Value = '6.94'
$ws.Range("E23").Value = '  -1.43%  '
$ws.Range("D24").Value = '82.06'
$ws.Range("E24").Value = '  +4.54%  '
$ws.Range("D25").Value = '11.94'
$ws.Range("E25").Value = '  -1.48%  '
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  -0.07%  '
$ws.Range("D27").Value = '2.69'
$ws.Range("E27").Value = '  -0.51%  '
$ws.Range("D28").Value = '7.71'
$ws.Range("E28").Value = '  -2.15%  '
$ws.Range("D29").Value = '0.998'
$ws.Range("E29").Value = '  -0.13%  '
$ws.Range("D30").Value = '1.91'
$ws.Range("E30").Value = '  +1.61%  '
$ws.Range("D31").Value = '25.71'
$ws.Range("E31").Value = '  -1.58%  '
$ws.Range("D32").Value = '1.11'
$ws.Range("E32").Value = '  -4.09%  '
$ws.Range("D33").Value = '5.62'
$ws.Range("E33").Value = '  +2.85%  '
$ws.Range("E34").Value = '  +0.81%  '
$ws.Range("E35").Value = '  -7.11%  '
$ws.Range("D36").Value = '5.84'
$ws.Range("E36").Value = '  -2.26%  '
$ws.Range("D37").Value = '3.138.17'
$ws.Range("E37").Value = '  -3.44%  '
$ws.Range("D38").Value = '436.48'
$ws.Range("E38").Value = '  -10.12%  '
$ws.Range("D39").Value = '0.0794'
$ws.Range("E39").Value = '  +0.18%  '
$ws.Range("D40").Value = '0.0382'
$ws.Range("E40").Value = '  -4.20%  '
$ws.Range("E41").Value = '  -0.45%  '
$ws.Range("D42").Value = '8.07'
$ws.Range("E42").Value = '  -0.73%  '
$ws.Range("D43").Value = '2.43'
$ws.Range("E43").Value = '  -5.66%  '
$ws.Range("D44").Value = '26.38'
$ws.Range("E44").Value = '  +4.30%  '
$ws.Range("D46").Value = '0.241'
$ws.Range("E46").Value = '  -4.04%  '
$ws.Range("E47").Value = '  -0.52%  '
$ws.Range("D48").Value = '1.94'
$ws.Range("E48").Value = '  -3.06%  '
$ws.Range("D49").Value = '115.32'
$ws.Range("E49").Value = '  -6.20%  '
$ws.Range("E50").Value = '  +4.38%  '
$ws.Range("D51").Value = '0.0₃0486'
$ws.Range("E51").Value = '  -6.73%  '

foreach ($r in $dRows) {
    $ws.Range("D$r").Style = "Normal"
}
